$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Current data occupies rows 2-22 (21 data rows), columns A-C.
# Capture the existing values for rows 2-15 (these are the rows that survive,
# shifting down by 6 to become rows 8-21).
$oldValues = @{}
for ($r = 2; $r -le 15; $r++) {
    $oldValues[$r] = @(
        $ws.Cells.Item($r, 1).Value(),
        $ws.Cells.Item($r, 2).Value(),
        $ws.Cells.Item($r, 3).Value()
    )
}

# Write the preserved rows into their new positions (8-21), working from the
# bottom up so we never overwrite a row before it has been read.
for ($r = 15; $r -ge 2; $r--) {
    $destRow = $r + 6
    $vals = $oldValues[$r]
    $ws.Cells.Item($destRow, 1).Value = $vals[0]
    $ws.Cells.Item($destRow, 2).Value = $vals[1]
    $ws.Cells.Item($destRow, 3).Value = $vals[2]
}

# New rows of accelerometer data inserted at the top (rows 2-7)
$newData = @(
    @(2.566667938232422, -3.378203916549682, 3.007539582252503),
    @(3.106618106365205, -3.249815458059311, 3.031012719869614),
    @(2.987140679359436, -3.142817544937134, 3.183629143238068),
    @(2.434188187122345, -3.181812554597855, 3.162444919347763),
    @(2.282221984863281, -3.265003252029419, 3.094355344772339),
    @(2.110153055191039, -3.195758980512619, 3.138975620269776)
)

for ($i = 0; $i -lt $newData.Count; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $newData[$i][0]
    $ws.Cells.Item($row, 2).Value = $newData[$i][1]
    $ws.Cells.Item($row, 3).Value = $newData[$i][2]
}

# Clear out the now-stale rows 22-28 left over from the original rows 16-22,
# since the sheet now only spans down to row 21.
$ws.Range("A22:C28").ClearContents()
